# Updated cryptos list (prices + 1h volume deltas refreshed from the source feed).
# Price/volume cells are stored as plain text (coinranking.com formats thousands
# separators with dots, e.g. "64.526.64"), so every write below is a text value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.526.64'
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").Value = '3.467.01'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '''574.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").Value = '''160.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '3.465.28'
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '''0.573'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.40%  '
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("D11").Value = '''0.121'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.74%  '
$ws.Range("D12").Value = '''0.436'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.98%  '
$ws.Range("D13").Value = '4.061.94'
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").Value = '''27.60'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.61%  '
$ws.Range("E16").Value = '  -8.74%  '
$ws.Range("D17").Value = '64.685.51'
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").Value = '3.522.18'
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("D19").Value = '''6.23'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.92%  '
$ws.Range("D20").Value = '''13.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.53%  '
$ws.Range("D21").Value = '''380.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("D22").Value = '''7.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.79%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '''1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '''72.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("E25").Value = '  -4.93%  '
$ws.Range("E26").Value = '  -1.26%  '
$ws.Range("D27").Value = '''9.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("E28").Value = '  +0.45%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '''6.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.94%  '
$ws.Range("E31").Value = '  -6.39%  '
$ws.Range("E32").Value = '  -1.49%  '
$ws.Range("D33").Value = '''23.34'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("D34").Value = '''7.04'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.52%  '
$ws.Range("D35").Value = '''1.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.50%  '
$ws.Range("D36").Value = '''160.89'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").Value = '''1.87'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.35%  '
$ws.Range("D38").Value = '''26.89'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.03%  '
$ws.Range("D39").Value = '''0.822'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.32%  '
$ws.Range("D40").Value = '''0.0747'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.43%  '
$ws.Range("D41").Value = '2.837.08'
$ws.Range("E41").Value = '  -2.69%  '
$ws.Range("E42").Value = '  -3.96%  '
$ws.Range("D43").Value = '''42.86'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.70%  '
$ws.Range("D44").Value = '''6.47'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.32%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '''0.0310'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.83%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '''25.84'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("D47").Value = '''2.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.76%  '
$ws.Range("D48").Value = '''335.50'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.22%  '
$ws.Range("E49").Value = '  -3.24%  '
$ws.Range("D50").Value = '''6.46'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.43%  '
$ws.Range("D51").Value = '''0.843'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.03%  '
